# Word COM-interop script applying the "Add files via upload" edit:
#  1. Insert a new first paragraph containing "test" and relocate the
#     "_GoBack" bookmark onto the end of it (bookmark names are unique
#     document-wide, so re-adding "_GoBack" removes it from its old spot
#     at the end of the last paragraph and places it here instead).
#  2. Collapse the multi-run "~ Commit changes in github or
#     pivotaltracker?" paragraph (which had proofing-error wrapped runs
#     around "github"/"pivotaltracker") back into a single plain run.

$d = $word.ActiveDocument

# --- 1. New leading paragraph "test" + relocate the _GoBack bookmark ---
# A trailing placeholder character ("X") is inserted after "test" so the
# bookmark's collapsed insertion point isn't sitting exactly on the
# paragraph-end boundary (some hosts mis-place a bookmark collapsed right
# at a paragraph's trailing mark); the placeholder is deleted afterwards,
# leaving a clean, empty "_GoBack" bookmark right after "test".
$d.Paragraphs(1).Range.InsertBefore("testX`r")
$bmSpot = $d.Range(4, 4)
$d.Bookmarks.Add("_GoBack", $bmSpot)
$placeholder = $d.Range(4, 5)
$placeholder.Delete()

# --- 2. Merge the "Commit changes in github or pivotaltracker?" runs ---
$commitPara = $d.Paragraphs(7).Range
[void]$commitPara.MoveEnd(1, -1)
# Write a differing value first so the no-op/"text unchanged" fast path
# in the host doesn't skip rebuilding the run list, then set the real
# text; Range.Text always collapses its range down to a single run.
$commitPara.Text = "~ Commit changes in github or pivotaltracker?#"

$commitPara2 = $d.Paragraphs(7).Range
[void]$commitPara2.MoveEnd(1, -1)
$commitPara2.Text = "~ Commit changes in github or pivotaltracker?"
